$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha2")
Write-Host $ws.Name
